$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => [new nombre_aides (C), new montant_total (E)]
$updates = @{
    3   = @(249334, 1036486441)
    91  = @(151202, 482926544)
    92  = @(409279, 1597133234)
    93  = @(209653, 1309991539)
    94  = @(94229,  918982287)
    95  = @(50801,  934231255)
    96  = @(17322,  797242656)
    104 = @(135300, 272650477)
    167 = @(12220,  105794259)
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Range("C$row").Value = $vals[0]
    $ws.Range("E$row").Value = $vals[1]
}
